# The edit swaps the distinguishing values between row 5 and row 6
# (columns A, I, Q, R, Y, AA, AC) while leaving every other column
# (which already hold identical data in both rows) untouched.
#
# Range.Copy(destination) is used instead of Value assignment so that
# text that looks like a date (e.g. "2022-09-12") is carried over as a
# literal string, exactly like the source cell, instead of being
# reinterpreted as a date serial number with a new number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "I", "Q", "R", "Y", "AA", "AC")
$tmpRow = 1000

foreach ($col in $cols) {
    $cell5 = $ws.Range("{0}5" -f $col)
    $cell6 = $ws.Range("{0}6" -f $col)
    $tmp   = $ws.Range("{0}{1}" -f $col, $tmpRow)

    # row5 -> tmp
    if ($cell5.Value2 -eq $null) {
        $tmp.Clear()
    } else {
        $cell5.Copy($tmp)
    }

    # row6 -> row5
    if ($cell6.Value2 -eq $null) {
        $cell5.Clear()
    } else {
        $cell6.Copy($cell5)
    }

    # tmp -> row6
    if ($tmp.Value2 -eq $null) {
        $cell6.Clear()
    } else {
        $tmp.Copy($cell6)
    }

    $tmp.Clear()
}
